$d = $word.ActiveDocument

$replacements = @(
    @("307×4=", "884×8="),
    @("402×7=", "313×5="),
    @("344×2=", "120×5="),
    @("618×7=", "556×6="),
    @("239×6=", "709×9="),
    @("185×4=", "288×6="),
    @("732×8=", "154×6="),
    @("647×8=", "358×4="),
    @("196×4=", "397×7="),
    @("667×3=", "121×3="),
    @("406×4=", "949×6="),
    @("229×8=", "997×4="),
    @("424×3=", "459×7="),
    @("541×9=", "493×6="),
    @("511×3=", "192×8="),
    @("305×6=", "552×6="),
    @("230×8=", "101×2="),
    @("964×2=", "848×6="),
    @("508×5=", "857×2="),
    @("624×9=", "287×3="),
    @("693×8=", "582×2="),
    @("844×4=", "238×9="),
    @("304×9=", "402×3="),
    @("895×8=", "874×4="),
    @("597×3=", "261×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
